$d = $word.ActiveDocument

$d.Content.Find.Execute("985×9=", $true, $false, $false, $false, $false, $true, 1, $false, "221×8=", 2) | Out-Null
$d.Content.Find.Execute("997×5=", $true, $false, $false, $false, $false, $true, 1, $false, "455×2=", 2) | Out-Null
$d.Content.Find.Execute("297×7=", $true, $false, $false, $false, $false, $true, 1, $false, "295×2=", 2) | Out-Null
$d.Content.Find.Execute("233×9=", $true, $false, $false, $false, $false, $true, 1, $false, "818×7=", 2) | Out-Null
$d.Content.Find.Execute("917×3=", $true, $false, $false, $false, $false, $true, 1, $false, "350×4=", 2) | Out-Null
$d.Content.Find.Execute("834×7=", $true, $false, $false, $false, $false, $true, 1, $false, "969×3=", 2) | Out-Null
$d.Content.Find.Execute("496×9=", $true, $false, $false, $false, $false, $true, 1, $false, "979×9=", 2) | Out-Null
$d.Content.Find.Execute("695×9=", $true, $false, $false, $false, $false, $true, 1, $false, "287×7=", 2) | Out-Null
$d.Content.Find.Execute("752×3=", $true, $false, $false, $false, $false, $true, 1, $false, "770×3=", 2) | Out-Null
$d.Content.Find.Execute("287×3=", $true, $false, $false, $false, $false, $true, 1, $false, "195×7=", 2) | Out-Null
$d.Content.Find.Execute("418×8=", $true, $false, $false, $false, $false, $true, 1, $false, "608×4=", 2) | Out-Null
$d.Content.Find.Execute("701×8=", $true, $false, $false, $false, $false, $true, 1, $false, "649×5=", 2) | Out-Null
$d.Content.Find.Execute("355×3=", $true, $false, $false, $false, $false, $true, 1, $false, "207×9=", 2) | Out-Null
$d.Content.Find.Execute("237×9=", $true, $false, $false, $false, $false, $true, 1, $false, "924×9=", 2) | Out-Null
$d.Content.Find.Execute("125×2=", $true, $false, $false, $false, $false, $true, 1, $false, "624×9=", 2) | Out-Null
$d.Content.Find.Execute("854×4=", $true, $false, $false, $false, $false, $true, 1, $false, "876×5=", 2) | Out-Null
$d.Content.Find.Execute("473×2=", $true, $false, $false, $false, $false, $true, 1, $false, "720×7=", 2) | Out-Null
$d.Content.Find.Execute("583×3=", $true, $false, $false, $false, $false, $true, 1, $false, "881×5=", 2) | Out-Null
$d.Content.Find.Execute("574×8=", $true, $false, $false, $false, $false, $true, 1, $false, "161×9=", 2) | Out-Null
$d.Content.Find.Execute("610×6=", $true, $false, $false, $false, $false, $true, 1, $false, "830×5=", 2) | Out-Null
$d.Content.Find.Execute("653×7=", $true, $false, $false, $false, $false, $true, 1, $false, "354×2=", 2) | Out-Null
$d.Content.Find.Execute("595×4=", $true, $false, $false, $false, $false, $true, 1, $false, "321×8=", 2) | Out-Null
$d.Content.Find.Execute("288×8=", $true, $false, $false, $false, $false, $true, 1, $false, "741×2=", 2) | Out-Null
$d.Content.Find.Execute("970×2=", $true, $false, $false, $false, $false, $true, 1, $false, "437×3=", 2) | Out-Null
$d.Content.Find.Execute("339×2=", $true, $false, $false, $false, $false, $true, 1, $false, "349×6=", 2) | Out-Null
